$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 4; existing rows 4..128 shift down to 5..129
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 (same constant columns as the rest of the
# dataset, plus the new date/volume/price data from the commit)
$ws.Cells.Item(4, 1).Value = 5
$ws.Cells.Item(4, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(4, 3).Value = "Maule"
$ws.Cells.Item(4, 4).Value = 44631
$ws.Cells.Item(4, 5).Value = 7
$ws.Cells.Item(4, 6).Value = 100112030
$ws.Cells.Item(4, 7).Value = "Poroto granado"
$ws.Cells.Item(4, 8).Value = "Sin especificar"
$ws.Cells.Item(4, 9).Value = "Primera"
$ws.Cells.Item(4, 10).Value = 300
$ws.Cells.Item(4, 11).Value = 20000
$ws.Cells.Item(4, 12).Value = 20000
$ws.Cells.Item(4, 13).Value = 20000
$ws.Cells.Item(4, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(4, 15).Value = "Región del Maule"
$ws.Cells.Item(4, 16).Value = 800
$ws.Cells.Item(4, 17).Value = 25
$ws.Cells.Item(4, 18).Value = "Hortaliza"

# Match the date style used by the rest of column D
$ws.Cells.Item(4, 4).NumberFormat = $ws.Cells.Item(5, 4).NumberFormat
